# code refactoring for watch list test cases
# Adds two new rows (TestCase_F15, TestCase_F16) to the "Test Cases" sheet,
# widens column B to fit the new Jira-id content, and moves the
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 16: TestCase_F15 -------------------------------------------------
# Copy formatting from the row right above (row 15) so the new row picks up
# the same fills/borders/wrap settings, then overwrite the values.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)

$ws.Range("A16").Value = "TestCase_F15"
$ws.Range("B16").Value = "OPQA-226"
$ws.Range("C16").Value = "Verify that users should be able to select from a list of suggested topics and check selected topic is presented in users type ahead"
$ws.Range("D16").Value = "Y"
$ws.Range("E16").Value = "SKIP"

# --- Row 17: TestCase_F16 --------------------------------------------------
# A17/D17/E17 reuse the plain "watermark" style (same as A15), B17 reuses the
# bordered/no-fill style (same as C14/C15), and C17 needs that same
# bordered/no-fill style PLUS word-wrap - toggling WrapText on a cell copied
# from C14 reuses the existing wrap+border style instead of minting a new one.
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").WrapText = $true

$ws.Range("D15").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("A15").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("A17").Value = "TestCase_F16"
$ws.Range("B17").Value = "OPQA-231,OPQA-1100"
$ws.Range("C17").Value = "Verify that Trending now section include articles and posts and able to navigate from tending now section and 
Verify that Maximum count on the trending list is 10"
$ws.Range("D17").Value = "Y"
$ws.Range("E17").Value = "PASS"

$ws.Rows.Item(17).RowHeight = 30

# --- Column widths ----------------------------------------------------------
# Column B now needs to be wide enough to fit "OPQA-231,OPQA-1100"; auto-fit
# it (this also naturally splits the old merged A:B column-width group since
# column A keeps its original width).
$ws.Columns("B:B").AutoFit()

# --- Selection / scroll position -------------------------------------------
$ws.Activate()
$ws.Range("D14").Select()
